{"js": "// Update the worksheet date and all the two-digit-division answers.\n// Each (oldText -> newText) pair below is unique within the document,\n// so we can safely locate-and-replace each one independently.\nconst replacements = [\n  [\"2023-09-18 Monday\", \"2023-09-19 Tuesday\"],\n  [\"96\u00f73=32, 0\", \"78\u00f76=13, 0\"],\n  [\"62\u00f79=6, 8\", \"20\u00f73=6, 2\"],\n  [\"83\u00f73=27, 2\", \"20\u00f77=2, 6\"],\n  [\"99\u00f73=33, 0\", \"50\u00f73=16, 2\"],\n  [\"76\u00f73=25, 1\", \"10\u00f77=1, 3\"],\n  [\"57\u00f75=11, 2\", \"38\u00f72=19, 0\"],\n  [\"88\u00f75=17, 3\", \"40\u00f75=8, 0\"],\n  [\"92\u00f75=18, 2\", \"98\u00f79=10, 8\"],\n  [\"46\u00f78=5, 6\", \"76\u00f79=8, 4\"],\n  [\"54\u00f72=27, 0\", \"27\u00f76=4, 3\"],\n  [\"33\u00f74=8, 1\", \"71\u00f75=14, 1\"],\n  [\"76\u00f74=19, 0\", \"73\u00f73=24, 1\"],\n  [\"34\u00f77=4, 6\", \"78\u00f79=8, 6\"],\n  [\"98\u00f75=19, 3\", \"62\u00f78=7, 6\"],\n  [\"23\u00f79=2, 5\", \"24\u00f74=6, 0\"],\n  [\"20\u00f78=2, 4\", \"10\u00f75=2, 0\"],\n  [\"18\u00f77=2, 4\", \"51\u00f78=6, 3\"],\n  [\"33\u00f75=6, 3\", \"72\u00f78=9, 0\"],\n  [\"98\u00f78=12, 2\", \"44\u00f76=7, 2\"],\n  [\"34\u00f73=11, 1\", \"65\u00f72=32, 1\"],\n  [\"11\u00f77=1, 4\", \"35\u00f77=5, 0\"],\n  [\"19\u00f72=9, 1\", \"95\u00f74=23, 3\"],\n  [\"87\u00f77=12, 3\", \"29\u00f73=9, 2\"],\n  [\"58\u00f78=7, 2\", \"10\u00f73=3, 1\"],\n  [\"37\u00f72=18, 1\", \"75\u00f77=10, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all the two-digit-division answers.\n# Each (old -> new) pair below is unique within the document, so a plain\n# Find/Replace for each exact string is sufficient and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-18 Monday\", \"2023-09-19 Tuesday\"),\n    @(\"96\u00f73=32, 0\", \"78\u00f76=13, 0\"),\n    @(\"62\u00f79=6, 8\", \"20\u00f73=6, 2\"),\n    @(\"83\u00f73=27, 2\", \"20\u00f77=2, 6\"),\n    @(\"99\u00f73=33, 0\", \"50\u00f73=16, 2\"),\n    @(\"76\u00f73=25, 1\", \"10\u00f77=1, 3\"),\n    @(\"57\u00f75=11, 2\", \"38\u00f72=19, 0\"),\n    @(\"88\u00f75=17, 3\", \"40\u00f75=8, 0\"),\n    @(\"92\u00f75=18, 2\", \"98\u00f79=10, 8\"),\n    @(\"46\u00f78=5, 6\", \"76\u00f79=8, 4\"),\n    @(\"54\u00f72=27, 0\", \"27\u00f76=4, 3\"),\n    @(\"33\u00f74=8, 1\", \"71\u00f75=14, 1\"),\n    @(\"76\u00f74=19, 0\", \"73\u00f73=24, 1\"),\n    @(\"34\u00f77=4, 6\", \"78\u00f79=8, 6\"),\n    @(\"98\u00f75=19, 3\", \"62\u00f78=7, 6\"),\n    @(\"23\u00f79=2, 5\", \"24\u00f74=6, 0\"),\n    @(\"20\u00f78=2, 4\", \"10\u00f75=2, 0\"),\n    @(\"18\u00f77=2, 4\", \"51\u00f78=6, 3\"),\n    @(\"33\u00f75=6, 3\", \"72\u00f78=9, 0\"),\n    @(\"98\u00f78=12, 2\", \"44\u00f76=7, 2\"),\n    @(\"34\u00f73=11, 1\", \"65\u00f72=32, 1\"),\n    @(\"11\u00f77=1, 4\", \"35\u00f77=5, 0\"),\n    @(\"19\u00f72=9, 1\", \"95\u00f74=23, 3\"),\n    @(\"87\u00f77=12, 3\", \"29\u00f73=9, 2\"),\n    @(\"58\u00f78=7, 2\", \"10\u00f73=3, 1\"),\n    @(\"37\u00f72=18, 1\", \"75\u00f77=10, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
